# The deck currently has the "Integral" design applied (its 12-colour
# DrawingML colour scheme lives in the theme part that the slide master /
# presentation point at). The edit being replayed here is the author
# switching the presentation's design back to the default "Office Theme"
# colour palette (Design tab -> Themes gallery -> Office).
#
# Applying a theme in PowerPoint rewrites the colour swatches of the
# currently-applied theme part in place; every slide inherits the new
# palette because they all resolve colours through
# Slide.ThemeColorScheme (backed by the slide master's theme).
#
# MsoThemeColorSchemeIndex order is:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5 Accent1, 6 Accent2, 7 Accent3, 8 Accent4, 9 Accent5, 10 Accent6,
#   11 Hyperlink, 12 FollowedHyperlink
# ColorFormat.RGB packs a colour the same way VBA's RGB() does:
#   value = R + G*256 + B*65536  (so e.g. 44546A -> 6968388)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$scheme = $s.ThemeColorScheme

$scheme.Item(1).RGB  = 0          # Dark 1    -> 000000
$scheme.Item(2).RGB  = 16777215   # Light 1   -> FFFFFF
$scheme.Item(3).RGB  = 6968388    # Dark 2    -> 44546A
$scheme.Item(4).RGB  = 15132391   # Light 2   -> E7E6E6
$scheme.Item(5).RGB  = 13998939   # Accent 1  -> 5B9BD5
$scheme.Item(6).RGB  = 3243501    # Accent 2  -> ED7D31
$scheme.Item(7).RGB  = 10855845   # Accent 3  -> A5A5A5
$scheme.Item(8).RGB  = 49407      # Accent 4  -> FFC000
$scheme.Item(9).RGB  = 12874308   # Accent 5  -> 4472C4
$scheme.Item(10).RGB = 4697456    # Accent 6  -> 70AD47
$scheme.Item(11).RGB = 12673797   # Hyperlink -> 0563C1
$scheme.Item(12).RGB = 7491477    # Followed Hyperlink -> 954F72
